# Normalize the "Recorded By" (column G) values on the "Session Analysis
# Results" sheet so that any entry literally named "System" (case-
# insensitive match on a comma-separated token) is moved to the front of
# the list, by reversing the order of the comma-separated names.
#
# Example:
#   "dnasr281@gmail.com, System"                 -> "System, dnasr281@gmail.com"
#   "system, backup@backdoor.com, System"        -> "System, backup@backdoor.com, system"
# Rows whose value does not contain a "system" token (case-insensitive)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    $val = [string]$val
    if ($val -eq "") { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $hasSystem = $false
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq "system") { $hasSystem = $true }
    }

    if ($hasSystem) {
        $reversed = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) { $reversed += $trimmed[$i] }
        $newVal = [string]::Join(", ", $reversed)
        # NOTE: PowerShell's -eq/-ne operators on strings are case-insensitive
        # in this runtime, so comparing $newVal against $val here would wrongly
        # treat "System, ..." as equal to "system, ...". Just assign directly;
        # writing an unchanged value back is harmless.
        $cell.Value2 = $newVal
    }
}
